$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Bitcoin", "BTC-USD", 719536829972.2048),
    @("Ethereum", "ETH-USD", 244206525520.8783),
    @("BNB", "BNB-USD", 37859945782.87637),
    @("XRP", "XRP-USD", 33915828443.44218),
    @("Solana", "SOL-USD", 26526553149.50745),
    @("Cardano", "ADA-USD", 13824714680.37449),
    @("Dogecoin", "DOGE-USD", 11189935570.81735),
    @("TRON", "TRX-USD", 9165671380.355772),
    @("Avalanche", "AVAX-USD", 8293800081.851548),
    @("Polygon", "MATIC-USD", 8283365709.063359),
    @("Toncoin", "TON-USD", 8010828924.011144),
    @("Chainlink", "LINK-USD", 7985886591.259349),
    @("Polkadot", "DOT-USD", 7060867585.842254),
    @("Wrapped Bitcoin", "WBTC-USD", 6018844475.097538),
    @("Litecoin", "LTC-USD", 5360223441.43322),
    @("Shiba Inu", "SHIB-USD", 5221139601.473591),
    @("Bitcoin Cash", "BCH-USD", 4596806846.195675),
    @("UNUS SED LEO", "LEO-USD", 3750487646.599951),
    @("Cosmos", "ATOM-USD", 3715654109.491059),
    @("OKB", "OKB-USD", 3523067013.188245),
    @("Stellar", "XLM-USD", 3377547916.808996),
    @("Kaspa", "KAS-USD", 3001318881.208693),
    @("Monero", "XMR-USD", 2946760695.065717),
    @("Ethereum Classic", "ETC-USD", 2876157252.276538),
    @("Cronos", "CRO-USD", 2513593377.987969)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
